# Corrected excel sheets for application fix issues

$wb = $excel.ActiveWorkbook

# --- Summary sheet: selection moves from B4 to A7 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A7").Select() | Out-Null

# --- Repayment schedule sheet: correct the % interest values and move
#     the selection from K7 to K2 (this sheet is no longer the active tab) ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Range("K2").Value = 5.56
$wsSchedule.Range("L2").Value = 5.56
$wsSchedule.Range("K2").Select() | Out-Null

# --- Transactions sheet: selection moves from A2:XFD4 (active A4) to F16 ---
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("F16").Select() | Out-Null

# --- NewLoanInput sheet: fix the product code text and make this sheet the
#     workbook's active tab (was "Repayment schedule") ---
$wsLoanInput = $wb.Worksheets.Item("NewLoanInput")
$wsLoanInput.Range("B2").Value = "2595-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-DISBURSE-FEE-%INT"
$wsLoanInput.Range("B2").Select() | Out-Null
